{"js": "// Amend the day in the \"Minutes\" document: Wednesday -> Tuesday\n// (row 6 of the minutes table: \"Group to meet either Wednesday or\n// Thursday to discuss final details before submitting document.\")\nconst results = context.document.body.search(\"Wednesday\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Wednesday\" in the document body.');\n}\n\n// Replace every match (there is exactly one in this document) in place,\n// preserving its original run formatting.\nresults.items.forEach((range) => {\n  range.insertText(\"Tuesday\", Word.InsertLocation.replace);\n});\nawait context.sync();\n", "ps1": "# Amend the day in the \"Minutes\" document: Wednesday -> Tuesday\n# (row 6 of the minutes table: \"Group to meet either Wednesday or\n# Thursday to discuss final details before submitting document.\")\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n#              Format, ReplaceWith, Replace)\n$find.Execute(\n    \"Wednesday\",  # FindText\n    $true,        # MatchCase\n    $true,        # MatchWholeWord\n    $false,       # MatchWildcards\n    $false,       # MatchSoundsLike\n    $false,       # MatchAllWordForms\n    $true,        # Forward\n    1,            # Wrap (wdFindContinue)\n    $false,       # Format\n    \"Tuesday\",    # ReplaceWith\n    2             # Replace (wdReplaceAll)\n) | Out-Null\n"}
